$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    # Forces the cell to hold TEXT even when $text looks like a number
    # (plain .Value assignment would otherwise coerce it to a Double).
    $c = $ws.Range($cellRef)
    $escaped = $text.Replace('"', '""')
    $c.Formula = '="' + $escaped + '"'
    $c.Copy() | Out-Null
    $c.PasteSpecial(-4163) | Out-Null
}
$excel.CutCopyMode = 0

$ws.Range("D2").Value = "57.856.74"
$ws.Range("E2").Value = "  +0.43%  "

$ws.Range("D3").Value = "3.134.60"
$ws.Range("E3").Value = "  +0.88%  "

$ws.Range("E4").Value = "  -0.07%  "

Set-TextValue "D5" "528.98"
$ws.Range("E5").Value = "  +1.14%  "

Set-TextValue "D6" "138.86"
$ws.Range("E6").Value = "  -1.25%  "

$ws.Range("E7").Value = "  -0.11%  "

$ws.Range("D8").Value = "3.132.55"
$ws.Range("E8").Value = "  +0.81%  "

Set-TextValue "D9" "0.447"
$ws.Range("E9").Value = "  +3.25%  "

Set-TextValue "D10" "7.19"
$ws.Range("E10").Value = "  -0.84%  "

$ws.Range("E11").Value = "  -1.20%  "

Set-TextValue "D12" "0.397"
$ws.Range("E12").Value = "  +3.05%  "

$ws.Range("D13").Value = "3.670.60"
$ws.Range("E13").Value = "  +0.47%  "

$ws.Range("E14").Value = "  +2.61%  "

Set-TextValue "D15" "25.52"
$ws.Range("E15").Value = "  -2.73%  "

$ws.Range("E16").Value = "  +0.44%  "

$ws.Range("D17").Value = "57.967.27"
$ws.Range("E17").Value = "  +0.46%  "

$ws.Range("D18").Value = "3.125.77"
$ws.Range("E18").Value = "  +0.21%  "

Set-TextValue "D19" "6.04"
$ws.Range("E19").Value = "  -0.98%  "

Set-TextValue "D20" "12.83"
$ws.Range("E20").Value = "  +0.08%  "

Set-TextValue "D21" "7.96"
$ws.Range("E21").Value = "  -1.27%  "

Set-TextValue "D22" "353.79"
$ws.Range("E22").Value = "  +5.20%  "

$ws.Range("E23").Value = "  +0.11%  "

Set-TextValue "D24" "68.57"
$ws.Range("E24").Value = "  +2.95%  "

Set-TextValue "D25" "0.507"
$ws.Range("E25").Value = "  -0.72%  "

Set-TextValue "D26" "0.170"
$ws.Range("E26").Value = "  +0.49%  "

$ws.Range("E27").Value = "  +0.00%  "

$ws.Range("D28").Value = "0.0₃0920"
$ws.Range("E28").Value = "  -0.29%  "

Set-TextValue "D29" "7.49"
$ws.Range("E29").Value = "  +3.96%  "

$ws.Range("E30").Value = "  +0.14%  "

Set-TextValue "D31" "6.22"
$ws.Range("E31").Value = "  -5.19%  "

Set-TextValue "D32" "1.89"
$ws.Range("E32").Value = "  +0.90%  "

Set-TextValue "D33" "21.16"
$ws.Range("E33").Value = "  +1.28%  "

$ws.Range("E34").Value = "  -0.74%  "

Set-TextValue "D35" "4.98"
$ws.Range("E35").Value = "  +6.99%  "

Set-TextValue "D36" "158.06"
$ws.Range("E36").Value = "  +1.25%  "

Set-TextValue "D37" "6.17"
$ws.Range("E37").Value = "  +1.27%  "

Set-TextValue "D38" "26.44"
$ws.Range("E38").Value = "  -1.94%  "

$ws.Range("E39").Value = "  -0.87%  "

$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D40" "0.0669"
$ws.Range("E40").Value = "  +0.92%  "

$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D41" "1.62"
$ws.Range("E41").Value = "  +4.76%  "

$ws.Range("E42").Value = "  +6.80%  "

Set-TextValue "D43" "0.704"
$ws.Range("E43").Value = "  +2.16%  "

$ws.Range("D44").Value = "3.171.44"
$ws.Range("E44").Value = "  +0.28%  "

$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D45" "36.62"
$ws.Range("E45").Value = "  -0.55%  "

$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D46" "0.0271"
$ws.Range("E46").Value = "  +5.05%  "

$ws.Range("E47").Value = "  -0.10%  "

$ws.Range("D48").Value = "2.325.68"
$ws.Range("E48").Value = "  +1.29%  "

Set-TextValue "D49" "0.973"
$ws.Range("E49").Value = "  -0.80%  "

$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D50" "20.53"
$ws.Range("E50").Value = "  -1.38%  "

$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D51" "6.04"
$ws.Range("E51").Value = "  +0.54%  "

$excel.CutCopyMode = 0
